$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    ,@(1, 'Page/Section', 'Layout', 'Type', 'Task', 'Priority')
    ,@(2, '*footer', 'All', 'programming', 'change social media icons to our relevant channels.', 'A')
    ,@(3, 'Home', 'All', 'programming', 'hero image slider autorotate every 6s', 'C')
    ,@(4, 'Home', 'All', 'programming', 'randomize initial image in rotator on page refresh', 'C')
    ,@(5, '*header', 'Desktop', 'programming', 'between widths 768 and 1035, the main nav breaks the layout - this can impact the logo carousel.', 'B')
    ,@(6, '*header', 'Mobile', 'programming', 'fix pulldown menu', 'A')
    ,@(7, '*footer', 'All', 'copy', 'Get TOS/PP from legal', 'A')
    ,@(8, 'Home', 'All', 'art', 'select hero imagery', 'A')
    ,@(9, 'Home', 'All', 'art', 'select imagery for business unit buttons (based on subsection hero imagery)', 'A')
    ,@(10, 'Home', 'All', 'copy', 'Revise homepage copy.', 'A')
    ,@(11, 'Pictures', 'All', 'art', 'select hero imagery (current image is FPO)', 'A')
    ,@(12, 'Pictures', 'All', 'copy', 'write description copy', 'A')
    ,@(13, 'Pictures', 'All', 'design', 'confirm presecnce of Select Film content', 'A')
    ,@(14, 'Pictures', 'All', 'art', 'select imagery for film posters', 'A')
    ,@(15, 'Press', 'All', 'art', 'select final hero imagery', 'B')
    ,@(16, 'Press', 'All', 'copy', 'write description copy', 'A')
    ,@(17, 'Press', 'All', 'copy', 'confirm presentation of projects', 'A')
    ,@(18, 'Adventures', 'All', 'copy', 'write description copy', 'A')
    ,@(19, 'Adventures', 'All', 'copy', 'confirm presentation of projects', 'A')
    ,@(20, 'Adventures', 'All', 'copy', 'confirm if Adventures will need a separate site', 'A')
    ,@(21, 'Promo', 'All', 'copy', 'write description copy', 'A')
    ,@(22, 'Promo', 'All', 'copy', 'confirm presentation of projects', 'A')
    ,@(23, 'Staffing', 'All', 'copy', 'confirm presentation of projects', 'A')
    ,@(24, 'Staffing', 'All', 'art', 'select hero imagery', 'A')
    ,@(25, 'Staffing', 'All', 'art', 'include staffing client logos', 'A')
    ,@(26, '*header', 'All', 'programming', 'the sticky header bounces when it transitions', 'D')
    ,@(27, 'Global', 'All', 'art', 'images need to be optimized - target: 4k resolution where possible; compression high (8).', 'B')
    ,@(28, 'Global', 'All', 'programming', 'have mobile sites load appropriately sized images ', 'C')
    ,@(29, 'Global', 'All', 'programming', 'implement Google Analytic (basic page views)', 'C')
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

# Remove the now-obsolete 30th row entirely so the used range / dimension shrinks to F29.
$ws.Rows.Item(30).Delete()

# Row heights: rows 6, 8, 12, 13, 14 lose their custom height (back to sheet default);
# rows 5, 9, 27 keep/gain the 30pt custom height used for wrapped, longer task text.
$ws.Rows.Item(6).EntireRow.AutoFit()
$ws.Rows.Item(8).EntireRow.AutoFit()
$ws.Rows.Item(12).EntireRow.AutoFit()
$ws.Rows.Item(13).EntireRow.AutoFit()
$ws.Rows.Item(14).EntireRow.AutoFit()

$ws.Rows.Item(5).RowHeight = 30
$ws.Rows.Item(9).RowHeight = 30
$ws.Rows.Item(27).RowHeight = 30

# Sheet view: drop the frozen topLeftCell scroll position and move the active
# selection to G5 (matches the author re-reviewing from the top of the sheet).
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("G5").Select()

# Data validations: rebuild the three list validations with their revised ranges
# now that the "Global" section exists and rows shifted.
$ws.Cells.Validation.Delete()

$sectionList = '"Home,Games,Pictures,Press,Adventures,Promo,Staffing,Agency,Technology,Capital,Contact,Careers,Our Story, *header,*footer,*metadata"'
$typeList = '"art,copy,programming,design"'
$sectionListWithGlobal = '"Home,Games,Pictures,Press,Adventures,Promo,Staffing,Agency,Technology,Capital,Contact,Careers,Our Story, Global, *header,*footer,*metadata"'

$dv1 = $ws.Range("B90:B97,A72:A97")
$dv1.Validation.Add(3, 1, 1, $sectionList)
$dv1.Validation.IgnoreBlank = $true
$dv1.Validation.InCellDropdown = $true
$dv1.Validation.ShowInput = $true
$dv1.Validation.ShowError = $true

$dv2 = $ws.Range("C2:C109")
$dv2.Validation.Add(3, 1, 1, $typeList)
$dv2.Validation.IgnoreBlank = $true
$dv2.Validation.InCellDropdown = $true
$dv2.Validation.ShowInput = $true
$dv2.Validation.ShowError = $true

$dv3 = $ws.Range("A2:A71")
$dv3.Validation.Add(3, 1, 1, $sectionListWithGlobal)
$dv3.Validation.IgnoreBlank = $true
$dv3.Validation.InCellDropdown = $true
$dv3.Validation.ShowInput = $true
$dv3.Validation.ShowError = $true

